{"js": "// Interface Control Document: bump revision to 1.03, log the change in the\n// revision-history table, and replace the \"N/A\" default-value placeholders\n// with \"Blank\" throughout the data-dictionary tables.\n\nconst body = context.document.body;\n\n// 1) Title page: \"Revision 1.02\" -> \"Revision 1.03\"\nconst revisionResults = body.search(\"Revision 1.02\", { matchCase: true });\nrevisionResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < revisionResults.items.length; i++) {\n  revisionResults.items[i].insertText(\"Revision 1.03\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Revision-history table: tidy up the 1.02 row's wrapped bullet so the\n// \" popular items\" / \" table\" runs become a single contiguous run\n// (\" popular items table\"), matching how Word recombines text when it is\n// retyped/reflowed.\nconst popularItemsResults = body.search(\" popular items table\", { matchCase: true });\npopularItemsResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < popularItemsResults.items.length; i++) {\n  popularItemsResults.items[i].insertText(\" popular items table\", \"Replace\");\n}\nawait context.sync();\n\n// 3) Data dictionary tables: replace every \"N/A\" default value with \"Blank\".\n// (Done before the revision-history row below is appended, so the literal\n// \"N/A\" mentioned in that row's description is left untouched.)\nconst naResults = body.search(\"N/A\", { matchCase: true });\nnaResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < naResults.items.length; i++) {\n  naResults.items[i].insertText(\"Blank\", \"Replace\");\n}\nawait context.sync();\n\n// 4) Append a new revision row documenting this change.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst revisionTable = tables.items[0];\nrevisionTable.addRows(\"End\", 1, [\n  [\"1.03\", \"6/27/17\", \"Jesse Cruse\", \"Changed \\u2018N/A\\u2019 default values to blank.\"]\n]);\nawait context.sync();\n", "ps1": "# Interface Control Document: bump revision to 1.03, log the change in the\n# revision-history table, and replace the \"N/A\" default-value placeholders\n# with \"Blank\" throughout the data-dictionary tables.\n\n$d = $word.ActiveDocument\n\n# 1) Title page: \"Revision 1.02\" -> \"Revision 1.03\"\n$d.Content.Find.Execute(\"Revision 1.02\", $false, $false, $false, $false, $false, $true, 1, $false, \"Revision 1.03\", 2) | Out-Null\n\n# 2) Revision-history table: tidy up the 1.02 row's wrapped bullet so the\n# \" popular items\" / \" table\" runs become a single contiguous run\n# (\" popular items table\"), matching how Word recombines text when it is\n# retyped/reflowed.\n$d.Content.Find.Execute(\" popular items table\", $true, $false, $false, $false, $false, $true, 1, $false, \" popular items table\", 2) | Out-Null\n\n# 3) Data dictionary tables: replace every \"N/A\" default value with \"Blank\".\n# (Done before the revision-history row below is appended, so the literal\n# \"N/A\" mentioned in that row's description is left untouched.)\n$d.Content.Find.Execute(\"N/A\", $false, $true, $false, $false, $false, $true, 1, $false, \"Blank\", 2) | Out-Null\n\n# 4) Append a new revision row documenting this change.\n$revisionTable = $d.Tables.Item(1)\n$newRow = $revisionTable.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"1.03\"\n$newRow.Cells.Item(2).Range.Text = \"6/27/17\"\n$newRow.Cells.Item(3).Range.Text = \"Jesse Cruse\"\n$newRow.Cells.Item(4).Range.Text = \"Changed \" + [char]0x2018 + \"N/A\" + [char]0x2019 + \" default values to blank.\"\n"}
